$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.776.32"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").Value = "2.320.03"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.64"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "82.94"
$ws.Range("E6").Value = "  -6.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0798"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.34"
$ws.Range("E11").Value = "  -7.84%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "2.680.75"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.31"
$ws.Range("E14").Value = "  -5.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.56"
$ws.Range("E15").Value = "  -6.94%  "
$ws.Range("D16").Value = "2.330.47"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.747"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "39.713.75"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").Value = "0.0₃0890"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.99"
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.72"
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.37"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.57"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  -6.28%  "
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.16"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.56"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.09"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.02"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.44"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0705"
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  -6.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0972"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("E39").Value = "  -9.08%  "
$ws.Range("E40").Value = "  -6.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("D42").Value = "1.964.23"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.26"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("E44").Value = "  -5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.16"
$ws.Range("E45").Value = "  -8.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.34"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -8.68%  "
$ws.Range("D48").Value = "2.541.99"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "91.55"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.24"
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.01"
$ws.Range("E51").Value = "  -4.87%  "
